$wb = $excel.ActiveWorkbook

# ALC row 7
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 2018.3334
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 2977.5
$ws.Range("K7").Value = 100
$ws.Range("L7").Value = 2977.5
$ws.Range("M7").Value = 12
$ws.Range("N7").Value = -3201.5

# ALC row 14
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 2018.3334
$ws.Range("I14").Value = 100
$ws.Range("J14").Value = 2977.5
$ws.Range("K14").Value = 100
$ws.Range("L14").Value = 2977.5
$ws.Range("M14").Value = 91
$ws.Range("N14").Value = -3359.5

# ALC row 48
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H48").Value = 10919
$ws.Range("I48").Value = 11103.8
$ws.Range("J48").Value = 9995
$ws.Range("K48").Value = 33311.39999999999
$ws.Range("L48").Value = 29985
$ws.Range("M48").Value = -33019.39999999999
$ws.Range("N48").Value = -30569

# ALC row 52
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 2409
$ws.Range("I52").Value = 899.6667
$ws.Range("J52").Value = 2975
$ws.Range("K52").Value = 2699.0001
$ws.Range("L52").Value = 8925
$ws.Range("M52").Value = -2539.0001
$ws.Range("N52").Value = -9245

# ALC row 56
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H56").Value = 10919
$ws.Range("I56").Value = 11103.8
$ws.Range("J56").Value = 9995
$ws.Range("K56").Value = 33311.39999999999
$ws.Range("L56").Value = 29985
$ws.Range("M56").Value = -32777.39999999999
$ws.Range("N56").Value = -31053

# ALC row 76
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 3056.389
$ws.Range("I76").Value = 2620.9092
$ws.Range("J76").Value = 3740.7144
$ws.Range("K76").Value = 2620.9092
$ws.Range("L76").Value = 3740.7144
$ws.Range("M76").Value = -2305.9092
$ws.Range("N76").Value = -4370.7144

# ALC row 79
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H79").Value = 3056.389
$ws.Range("I79").Value = 2620.9092
$ws.Range("J79").Value = 3740.7144
$ws.Range("K79").Value = 2620.9092
$ws.Range("L79").Value = 3740.7144
$ws.Range("M79").Value = -1528.9092
$ws.Range("N79").Value = -5924.7144

# ALC row 82
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H82").Value = 1025.2
$ws.Range("I82").Value = 1025.2
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3075.6
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2669.6
$ws.Range("N82").ClearContents()

# ALC row 85
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H85").Value = 1025.2
$ws.Range("I85").Value = 1025.2
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3075.6
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1671.6
$ws.Range("N85").ClearContents()

# ALC row 135
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 474.70587
$ws.Range("I135").Value = 497.5
$ws.Range("J135").Value = 420
$ws.Range("K135").Value = 4477.5
$ws.Range("L135").Value = 3780
$ws.Range("M135").Value = -1942.5
$ws.Range("N135").Value = -8850

# ARM row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 20001664
$ws.Range("I74").Value = 27779008
$ws.Range("J74").Value = 2775.4285
$ws.Range("K74").Value = 27779008
$ws.Range("L74").Value = 2775.4285
$ws.Range("M74").Value = -27778134
$ws.Range("N74").Value = -4523.4285

# ARM row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 20001664
$ws.Range("I77").Value = 27779008
$ws.Range("J77").Value = 2775.4285
$ws.Range("K77").Value = 138895040
$ws.Range("L77").Value = 13877.1425
$ws.Range("M77").Value = -138890672
$ws.Range("N77").Value = -22613.1425

# ARM row 141
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H141").Value = 75776.336
$ws.Range("J141").Value = 75776.336
$ws.Range("L141").Value = 75776.336
$ws.Range("N141").Value = -86136.336

# CRP row 10
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 1007
$ws.Range("I10").Value = 1007
$ws.Range("K10").Value = 1007
$ws.Range("M10").Value = -868

# CRP row 69
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H69").Value = 30600
$ws.Range("I69").Value = 14000
$ws.Range("J69").Value = 38900
$ws.Range("K69").Value = 14000
$ws.Range("L69").Value = 38900
$ws.Range("M69").Value = -13251
$ws.Range("N69").Value = -40398

# CRP row 72
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H72").Value = 30600
$ws.Range("I72").Value = 14000
$ws.Range("J72").Value = 38900
$ws.Range("K72").Value = 42000
$ws.Range("L72").Value = 116700
$ws.Range("M72").Value = -38256
$ws.Range("N72").Value = -124188

# CRP row 140
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 62645
$ws.Range("J140").Value = 62645
$ws.Range("L140").Value = 62645
$ws.Range("N140").Value = -73005

# CUL row 70
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2218.8
$ws.Range("I70").Value = 1548.5
$ws.Range("J70").Value = 4900
$ws.Range("K70").Value = 4645.5
$ws.Range("L70").Value = 14700
$ws.Range("M70").Value = -4330.5
$ws.Range("N70").Value = -15330

# CUL row 73
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H73").Value = 2218.8
$ws.Range("I73").Value = 1548.5
$ws.Range("J73").Value = 4900
$ws.Range("K73").Value = 4645.5
$ws.Range("L73").Value = 14700
$ws.Range("M73").Value = -3553.5
$ws.Range("N73").Value = -16884

# CUL row 98
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 300
$ws.Range("I98").Value = 250
$ws.Range("J98").Value = 400
$ws.Range("K98").Value = 750
$ws.Range("L98").Value = 1200
$ws.Range("M98").Value = 748
$ws.Range("N98").Value = -4196

# CUL row 99
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 1876.75
$ws.Range("I99").Value = 1202.8
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 3608.4
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = -1362.4
$ws.Range("N99").Value = -13492

# CUL row 100
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 6670.3335
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 6670.3335
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 20011.0005
$ws.Range("N100").Value = -21633.0005
$ws.Range("M100").ClearContents()

# CUL row 103
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 1862.6666
$ws.Range("I103").Value = 529.44446
$ws.Range("J103").Value = 2529.2778
$ws.Range("K103").Value = 1588.33338
$ws.Range("L103").Value = 7587.8334
$ws.Range("M103").Value = -709.33338
$ws.Range("N103").Value = -9345.8334

# CUL row 106
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H106").Value = 5669.8887
$ws.Range("J106").Value = 5669.8887
$ws.Range("L106").Value = 17009.6661
$ws.Range("N106").Value = -18901.6661

# CUL row 109
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 5556855.5
$ws.Range("I109").Value = 1950
$ws.Range("K109").Value = 5850
$ws.Range("M109").Value = -4810

# CUL row 125
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 4348.75
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 4348.75
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 13046.25
$ws.Range("N125").Value = -22886.25
$ws.Range("M125").ClearContents()

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 827.5469000000001
$ws.Range("I131").Value = 288.6154
$ws.Range("J131").Value = 964.9216
$ws.Range("K131").Value = 865.8462000000001
$ws.Range("L131").Value = 2894.7648
$ws.Range("M131").Value = 4174.1538
$ws.Range("N131").Value = -12974.7648

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 85508.75
$ws.Range("I80").Value = 2140
$ws.Range("J80").Value = 145057.86
$ws.Range("K80").Value = 2140
$ws.Range("L80").Value = 145057.86
$ws.Range("M80").Value = -1142
$ws.Range("N80").Value = -147053.86

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 85508.75
$ws.Range("I83").Value = 2140
$ws.Range("J83").Value = 145057.86
$ws.Range("K83").Value = 10700
$ws.Range("L83").Value = 725289.2999999999
$ws.Range("M83").Value = -5708
$ws.Range("N83").Value = -735273.2999999999

# GSM row 121
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2475.3333
$ws.Range("I122").Value = 3279.2
$ws.Range("J122").Value = 1901.1428
$ws.Range("K122").Value = 9837.599999999999
$ws.Range("L122").Value = 5703.428400000001
$ws.Range("M122").Value = -7387.599999999999
$ws.Range("N122").Value = -10603.4284
